# Weekly update: insert a new Coco (Mercado Mayorista Lo Valledor de Santiago)
# price record above the existing row 54, shifting rows 54-58 down to 55-59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54 (pushes old rows 54:58 down to 55:59).
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Cells.Item(54, 1).Value = 6
$ws.Cells.Item(54, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(54, 3).Value = "Metropolitana"
$ws.Cells.Item(54, 4).Value = 44474
$ws.Cells.Item(54, 5).Value = 13
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100108
$ws.Cells.Item(54, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(54, 9).Value = 100108007
$ws.Cells.Item(54, 10).Value = "Coco"
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 250
$ws.Cells.Item(54, 14).Value = 23000
$ws.Cells.Item(54, 15).Value = 25000
$ws.Cells.Item(54, 16).Value = 24600
$ws.Cells.Item(54, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(54, 18).Value = "Perú"
$ws.Cells.Item(54, 19).Value = 1230
$ws.Cells.Item(54, 20).Value = 20
